$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows below the existing "S 1" week block (rows 8-11) so the
# schedule grid grows to hold a second week, plus one trailing blank row.
$ws.Rows.Item(12).Resize(5).Insert()

# Duplicate the week's schedule (D8:I11) down into the newly inserted rows
# (D12:I15) as the template for the second week.
$src = $ws.Range("D8:I11")
$dst = $ws.Range("D12:I15")
$src.Copy()
$dst.PasteSpecial(-4104)

# The bulk paste above only reliably carries the per-row number format for
# columns D:H; column I's "T.I"/"R.O" Arial style (I8, I10) needs a
# dedicated formats-only touch-up onto its new-week counterparts (I12, I14).
$ws.Range("I8").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I10").Copy()
$ws.Range("I14").PasteSpecial(-4122)

# Relabel "S 1" -> "S 11" (first week) and the new block -> "S 12" (second week).
$ws.Range("D8:D11").Value2 = "S 11"
$ws.Range("D12:D15").Value2 = "S 12"

# Give the Wednesday "C#" module cells (I11, I15) the real-time highlight
# font: same Arial face as the other module cells, but automatic/theme text
# color instead of the muted grey used elsewhere.
$ws.Range("I8").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Font.ThemeColor = 1

$ws.Range("I8").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I15").Font.ThemeColor = 1

# Mirror column I's formatting one row further so the trailing blank row
# (16) is ready for new entries, matching G16/H16.
$ws.Range("I10").Copy()
$ws.Range("I16").PasteSpecial(-4122)
$ws.Range("I16").ClearContents()

$ws.Application.CutCopyMode = $false

$ws.Range("J19").Select()
